$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Property1" to "DataNode" (unify DataNode/DataTable/Entity naming)
$ws.Name = "DataNode"

# Update the active selection in the frozen (bottomLeft) pane from K9 to O40
$ws.Range("O40").Select()
